$wb = $excel.ActiveWorkbook

# F-column ("想去人数" / want-to-go count) updates, applied identically
# to the "展览" and "全部类型" sheets (they mirror the same event data).
$updates = @{
    "F2" = 1193
    "F3" = 977
    "F4" = 294
    "F8" = 2428
    "F9" = 7936
    "F10" = 942
    "F11" = 478
    "F12" = 420
    "F13" = 184
    "F14" = 446
    "F15" = 10
    "F17" = 8179
    "F18" = 327
    "F19" = 1413
    "F23" = 195
    "F24" = 346
    "F27" = 26
    "F28" = 117
    "F31" = 1170
    "F32" = 31
    "F33" = 59
    "F34" = 105
    "F35" = 71
    "F36" = 90
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}

Write-Output "Updated F-column counts on 展览 and 全部类型 sheets"
